$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 334
$ws.Range("I4").Value = 167.5
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 167.5
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -53.5
$ws.Range("N4").Value = -1228
$ws.Range("N29").ClearContents()
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H43").Value = 400
$ws.Range("I43").Value = 400
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 400
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -331
$ws.Range("H80").Value = 7167444.5
$ws.Range("I80").Value = 711.9286
$ws.Range("J80").Value = 12184157
$ws.Range("K80").Value = 2135.7858
$ws.Range("L80").Value = 36552471
$ws.Range("M80").Value = -1137.7858
$ws.Range("N80").Value = -36554467
$ws.Range("H83").Value = 7167444.5
$ws.Range("I83").Value = 711.9286
$ws.Range("J83").Value = 12184157
$ws.Range("K83").Value = 6407.3574
$ws.Range("L83").Value = 109657413
$ws.Range("M83").Value = -1415.3574
$ws.Range("N83").Value = -109667397
$ws.Range("H129").Value = 164945.45
$ws.Range("J129").Value = 186283.62
$ws.Range("L129").Value = 558850.86
$ws.Range("N129").Value = -568850.86
$ws.Range("N130").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H132").Value = 5094.5
$ws.Range("I132").Value = 5367.467
$ws.Range("K132").Value = 16102.401
$ws.Range("M132").Value = -13572.401
$ws.Range("H137").Value = 18659.154
$ws.Range("I137").Value = 1236.5952
$ws.Range("J137").Value = 64393.375
$ws.Range("K137").Value = 3709.7856
$ws.Range("L137").Value = 193180.125
$ws.Range("M137").Value = -1159.7856
$ws.Range("N137").Value = -198280.125
$ws.Range("H138").Value = 2413.5076
$ws.Range("I138").Value = 3496
$ws.Range("J138").Value = 2266.7288
$ws.Range("K138").Value = 10488
$ws.Range("L138").Value = 6800.1864
$ws.Range("M138").Value = -5348
$ws.Range("N138").Value = -17080.1864
$ws.Range("H141").Value = 1477.3103
$ws.Range("I141").Value = 1001.7917
$ws.Range("J141").Value = 3759.8
$ws.Range("K141").Value = 3005.3751
$ws.Range("L141").Value = 11279.4
$ws.Range("M141").Value = 2174.6249
$ws.Range("N141").Value = -21639.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35993.43
$ws.Range("I32").Value = 41531.9
$ws.Range("J32").Value = 2762.6
$ws.Range("K32").Value = 41531.9
$ws.Range("L32").Value = 2762.6
$ws.Range("M32").Value = -41244.9
$ws.Range("N32").Value = -3336.6
$ws.Range("H45").Value = 2615.5881
$ws.Range("I45").Value = 2341.3333
$ws.Range("K45").Value = 2341.3333
$ws.Range("M45").Value = -1964.3333
$ws.Range("N74").Value = -2747.75
$ws.Range("H74").Value = 41667584
$ws.Range("I74").Value = 50000900
$ws.Range("J74").Value = 999.75
$ws.Range("K74").Value = 50000900
$ws.Range("L74").Value = 999.75
$ws.Range("M74").Value = -50000026
$ws.Range("N77").Value = -13734.75
$ws.Range("H77").Value = 41667584
$ws.Range("I77").Value = 50000900
$ws.Range("J77").Value = 999.75
$ws.Range("K77").Value = 250004500
$ws.Range("L77").Value = 4998.75
$ws.Range("M77").Value = -250000132
$ws.Range("H132").Value = 30279.697
$ws.Range("I132").Value = 1725.5682
$ws.Range("J132").Value = 169877.67
$ws.Range("K132").Value = 5176.7046
$ws.Range("L132").Value = 509633.01
$ws.Range("M132").Value = -2646.7046
$ws.Range("N132").Value = -514693.01
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 548.73914
$ws.Range("I94").Value = 526.1
$ws.Range("J94").Value = 699.6667
$ws.Range("K94").Value = 526.1
$ws.Range("L94").Value = 699.6667
$ws.Range("M94").Value = -75.10000000000002
$ws.Range("N94").Value = -1601.6667
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 494
$ws.Range("I7").Value = 494
$ws.Range("K7").Value = 494
$ws.Range("M7").Value = -381
$ws.Range("H52").Value = 39949.5
$ws.Range("J52").Value = 39949.5
$ws.Range("L52").Value = 39949.5
$ws.Range("N52").Value = -40537.5
$ws.Range("H58").Value = 14050.632
$ws.Range("I58").Value = 937.37933
$ws.Range("J58").Value = 56304.445
$ws.Range("K58").Value = 937.37933
$ws.Range("L58").Value = 56304.445
$ws.Range("M58").Value = -734.37933
$ws.Range("N58").Value = -56710.445
$ws.Range("N70").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("H132").Value = 23272.92
$ws.Range("I132").Value = 28567.525
$ws.Range("K132").Value = 85702.57500000001
$ws.Range("M132").Value = -83172.57500000001
$ws.Range("H136").Value = 14050.632
$ws.Range("I136").Value = 937.37933
$ws.Range("J136").Value = 56304.445
$ws.Range("K136").Value = 2812.13799
$ws.Range("L136").Value = 168913.335
$ws.Range("M136").Value = -262.1379900000002
$ws.Range("N136").Value = -174013.335
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4285940
$ws.Range("I4").Value = 344.75
$ws.Range("J4").Value = 10000067
$ws.Range("K4").Value = 1034.25
$ws.Range("L4").Value = 30000201
$ws.Range("M4").Value = -922.25
$ws.Range("N4").Value = -30000425
$ws.Range("M51").ClearContents()
$ws.Range("H51").Value = 3348.75
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("H131").Value = 751.26
$ws.Range("J131").Value = 803.32587
$ws.Range("L131").Value = 2409.97761
$ws.Range("N131").Value = -12489.97761
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N39").ClearContents()
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("H58").Value = 12505688
$ws.Range("I58").Value = 4000
$ws.Range("J58").Value = 16672917
$ws.Range("K58").Value = 4000
$ws.Range("L58").Value = 16672917
$ws.Range("M58").Value = -3723
$ws.Range("N58").Value = -16673471
$ws.Range("H102").Value = 20835616
$ws.Range("I102").Value = 27780326
$ws.Range("K102").Value = 27780326
$ws.Range("M102").Value = -27778704
$ws.Range("H113").Value = 2264.7144
$ws.Range("I113").Value = 2058.7144
$ws.Range("J113").Value = 2470.7144
$ws.Range("K113").Value = 2058.7144
$ws.Range("L113").Value = 2470.7144
$ws.Range("M113").Value = 111.2856000000002
$ws.Range("N113").Value = -6810.7144
$ws.Range("H136").Value = 19556.25
$ws.Range("J136").Value = 19556.25
$ws.Range("L136").Value = 58668.75
$ws.Range("N136").Value = -63768.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4276.2
$ws.Range("J22").Value = 10000
$ws.Range("L22").Value = 10000
$ws.Range("N22").Value = -10590
$ws.Range("H27").Value = 4276.2
$ws.Range("J27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("N27").Value = -10214
$ws.Range("N97").ClearContents()
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("H100").Value = 1676.0769
$ws.Range("I100").Value = 1257.1428
$ws.Range("J100").Value = 2164.8333
$ws.Range("K100").Value = 1257.1428
$ws.Range("L100").Value = 2164.8333
$ws.Range("M100").Value = -716.1428000000001
$ws.Range("N100").Value = -3246.8333
$ws.Range("H132").Value = 1492.5333
$ws.Range("I132").Value = 1106.8462
$ws.Range("K132").Value = 3320.5386
$ws.Range("M132").Value = -790.5385999999999
$ws.Range("H136").Value = 30911.234
$ws.Range("J136").Value = 2985
$ws.Range("L136").Value = 8955
$ws.Range("N136").Value = -14055
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N3").Value = -50000228
$ws.Range("H3").Value = 50000000
$ws.Range("J3").Value = 50000000
$ws.Range("L3").Value = 50000000
$ws.Range("M64").Value = -9552
$ws.Range("H64").Value = 19266.666
$ws.Range("I64").Value = 9800
$ws.Range("J64").Value = 24000
$ws.Range("K64").Value = 9800
$ws.Range("L64").Value = 24000
$ws.Range("N64").Value = -24496
$ws.Range("M67").Value = -8942
$ws.Range("H67").Value = 19266.666
$ws.Range("I67").Value = 9800
$ws.Range("J67").Value = 24000
$ws.Range("K67").Value = 9800
$ws.Range("L67").Value = 24000
$ws.Range("N67").Value = -25716
$ws.Range("H100").Value = 538.4
$ws.Range("I100").Value = 548
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 1096
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -555
$ws.Range("N100").Value = -2082
$ws.Range("H107").Value = 2841500.5
$ws.Range("I107").Value = 731.9
$ws.Range("K107").Value = 2195.7
$ws.Range("M107").Value = -275.6999999999998
$ws.Range("H122").Value = 1253.2354
$ws.Range("I122").Value = 1100.8334
$ws.Range("K122").Value = 3302.5002
$ws.Range("M122").Value = -852.5001999999999
$ws.Range("H136").Value = 27779180
$ws.Range("I136").Value = 30304242
$ws.Range("K136").Value = 90912726
$ws.Range("M136").Value = -90910176
